$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 12) -----------------------------------------------
# Add two more header columns ("Fin" and "Comentario" move from G/H to I/J)
# and insert a brand-new "Fecha" column. Clone the existing header format
# (bold/centered/bordered style) onto the two freshly used columns before
# the source cells get overwritten.
$ws.Range("B12").Copy()
$ws.Range("I12:J12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Shift the existing header labels two columns to the right: B->D, C->E,
# D->F, E->G, F->H, G->I, H->J (write back-to-front so nothing is clobbered).
$ws.Range("J12").Value = "Comentario"
$ws.Range("I12").Value = "Fin"
$ws.Range("H12").Value = "Inicio"
$ws.Range("G12").Value = "Tiempo"
$ws.Range("F12").Value = "Dinero generado"
$ws.Range("E12").Value = "Consola"
$ws.Range("D12").Value = "Empleado"

# New "Fecha" header fills the vacated C12, and B12 is now unused.
$ws.Range("C12").Value = "Fecha"
$ws.Range("B12").Clear()

# --- Column widths -------------------------------------------------------
# The two newly used columns (I, J) pick up an auto-fit-style width.
$ws.Columns.Item(9).ColumnWidth = 14.7
$ws.Columns.Item(10).ColumnWidth = 14.7

# --- Selection -------------------------------------------------------
$null = $ws.Range("L23").Select()
